$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7527
$ws1.Range("F7").Value = 928
$ws1.Range("F9").Value = 784
$ws1.Range("F10").Value = 575
$ws1.Range("F11").Value = 82
$ws1.Range("F12").Value = 65
$ws1.Range("F13").Value = 375
$ws1.Range("F14").Value = 855
$ws1.Range("F15").Value = 3020
$ws1.Range("F16").Value = 177
$ws1.Range("F17").Value = 73
$ws1.Range("F18").Value = 703
$ws1.Range("F19").Value = 748
$ws1.Range("F20").Value = 43
$ws1.Range("F23").Value = 188
$ws1.Range("F24").Value = 207
$ws1.Range("F25").Value = 226
$ws1.Range("F26").Value = 249
$ws1.Range("F27").Value = 122
$ws1.Range("F28").Value = 90
$ws1.Range("F29").Value = 228
$ws1.Range("F32").Value = 385
$ws1.Range("F33").Value = 431
$ws1.Range("F37").Value = 74

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 7527
$ws4.Range("F10").Value = 928
$ws4.Range("F12").Value = 784
$ws4.Range("F13").Value = 575
$ws4.Range("F14").Value = 82
$ws4.Range("F15").Value = 65
$ws4.Range("F16").Value = 375
$ws4.Range("F17").Value = 855
$ws4.Range("F19").Value = 3020
$ws4.Range("F20").Value = 177
$ws4.Range("F21").Value = 73
$ws4.Range("F23").Value = 703
$ws4.Range("F24").Value = 748
$ws4.Range("F26").Value = 43
$ws4.Range("F29").Value = 188
$ws4.Range("F30").Value = 207
$ws4.Range("F31").Value = 226
$ws4.Range("F32").Value = 249
$ws4.Range("F33").Value = 122
$ws4.Range("F34").Value = 90
$ws4.Range("F35").Value = 228
$ws4.Range("F38").Value = 385
$ws4.Range("F39").Value = 431
$ws4.Range("F43").Value = 74
